# chore: update Sheets via scheduled runner
# Refresh cached market-board derived figures (currentAveragePrice*, Leve
# price/profit columns H:N) across the per-job leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("N69").Value = 0

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = 0
$ws.Range("N72").Value = 0

$ws.Range("H106").Value = 8821789
$ws.Range("I106").Value = 15433694
$ws.Range("J106").Value = 5916.5
$ws.Range("K106").Value = 15433694
$ws.Range("L106").Value = 5916.5
$ws.Range("M106").Value = -15433063
$ws.Range("N106").Value = -7178.5

$ws.Range("H132").Value = 3788.8838
$ws.Range("I132").Value = 3967.139
$ws.Range("K132").Value = 11901.417
$ws.Range("M132").Value = -9371.417000000001

$ws.Range("H137").Value = 420414.2
$ws.Range("I137").Value = 738304.4399999999
$ws.Range("J137").Value = 15826.637
$ws.Range("K137").Value = 2214913.32
$ws.Range("L137").Value = 47479.911
$ws.Range("M137").Value = -2212363.32
$ws.Range("N137").Value = -52579.911

$ws.Range("H138").Value = 6064.7754
$ws.Range("J138").Value = 6789.9287
$ws.Range("L138").Value = 20369.7861
$ws.Range("N138").Value = -30649.7861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 168562.61
$ws.Range("I45").Value = 240829.44
$ws.Range("J45").Value = 5962.25
$ws.Range("K45").Value = 240829.44
$ws.Range("L45").Value = 5962.25
$ws.Range("M45").Value = -240452.44
$ws.Range("N45").Value = -6716.25

$ws.Range("H74").Value = 1731.5491
$ws.Range("I74").Value = 1427.3864
$ws.Range("K74").Value = 1427.3864
$ws.Range("M74").Value = -553.3864000000001

$ws.Range("H77").Value = 1731.5491
$ws.Range("I77").Value = 1427.3864
$ws.Range("K77").Value = 7136.932000000001
$ws.Range("M77").Value = -2768.932000000001

$ws.Range("H122").Value = 780755.5600000001
$ws.Range("I122").Value = 2335
$ws.Range("J122").Value = 1753781.2
$ws.Range("K122").Value = 7005
$ws.Range("L122").Value = 5261343.6
$ws.Range("M122").Value = -4555
$ws.Range("N122").Value = -5266243.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9010.583000000001
$ws.Range("I86").Value = 10042
$ws.Range("J86").Value = 3853.5
$ws.Range("K86").Value = 10042
$ws.Range("L86").Value = 3853.5
$ws.Range("M86").Value = -8919
$ws.Range("N86").Value = -6099.5

$ws.Range("H89").Value = 9010.583000000001
$ws.Range("I89").Value = 10042
$ws.Range("J89").Value = 3853.5
$ws.Range("K89").Value = 50210
$ws.Range("L89").Value = 19267.5
$ws.Range("M89").Value = -44594
$ws.Range("N89").Value = -30499.5

$ws.Range("H105").Value = 36032.863
$ws.Range("I105").Value = 41201.68
$ws.Range("K105").Value = 41201.68
$ws.Range("M105").Value = -39454.68

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2587.5
$ws.Range("I31").Value = 1706
$ws.Range("J31").Value = 3257.44
$ws.Range("K31").Value = 1706
$ws.Range("L31").Value = 3257.44
$ws.Range("M31").Value = -1411
$ws.Range("N31").Value = -3847.44

$ws.Range("H34").Value = 2587.5
$ws.Range("I34").Value = 1706
$ws.Range("J34").Value = 3257.44
$ws.Range("K34").Value = 1706
$ws.Range("L34").Value = 3257.44
$ws.Range("M34").Value = -1504
$ws.Range("N34").Value = -3661.44

$ws.Range("H132").Value = 24982.65
$ws.Range("I132").Value = 6782.1797
$ws.Range("J132").Value = 202437.25
$ws.Range("K132").Value = 20346.5391
$ws.Range("L132").Value = 607311.75
$ws.Range("M132").Value = -17816.5391
$ws.Range("N132").Value = -612371.75

$ws.Range("H134").Value = 5220016.5
$ws.Range("I134").Value = 5694363.5
$ws.Range("J134").Value = 2200
$ws.Range("K134").Value = 17083090.5
$ws.Range("L134").Value = 6600
$ws.Range("M134").Value = -17080555.5
$ws.Range("N134").Value = -11670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 401063.97
$ws.Range("I5").Value = 1363.4615
$ws.Range("J5").Value = 834072.8
$ws.Range("K5").Value = 4090.3845
$ws.Range("L5").Value = 2502218.4
$ws.Range("M5").Value = -3978.3845
$ws.Range("N5").Value = -2502442.4

$ws.Range("H12").Value = 13.818182
$ws.Range("I12").Value = 2.6
$ws.Range("J12").Value = 23.166666
$ws.Range("K12").Value = 7.800000000000001
$ws.Range("L12").Value = 69.49999800000001
$ws.Range("M12").Value = 165.2
$ws.Range("N12").Value = -415.499998

$ws.Range("H23").Value = 15151951
$ws.Range("I23").Value = 384.6
$ws.Range("J23").Value = 27778258
$ws.Range("K23").Value = 1153.8
$ws.Range("L23").Value = 83334774
$ws.Range("M23").Value = -918.8000000000002
$ws.Range("N23").Value = -83335244

$ws.Range("H68").Value = 8759.691999999999
$ws.Range("I68").Value = 3064.889
$ws.Range("J68").Value = 11774.588
$ws.Range("K68").Value = 9194.667000000001
$ws.Range("L68").Value = 35323.764
$ws.Range("M68").Value = -8383.667000000001
$ws.Range("N68").Value = -36945.764

$ws.Range("H71").Value = 8759.691999999999
$ws.Range("I71").Value = 3064.889
$ws.Range("J71").Value = 11774.588
$ws.Range("K71").Value = 27584.001
$ws.Range("L71").Value = 105971.292
$ws.Range("M71").Value = -23528.001
$ws.Range("N71").Value = -114083.292

$ws.Range("H107").Value = 7000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 7000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 21000
$ws.Range("M107").Value = 21000
$ws.Range("N107").Value = -24840

$ws.Range("H135").Value = 401063.97
$ws.Range("I135").Value = 1363.4615
$ws.Range("J135").Value = 834072.8
$ws.Range("K135").Value = 12271.1535
$ws.Range("L135").Value = 7506655.2
$ws.Range("M135").Value = -9736.153499999999
$ws.Range("N135").Value = -7511725.2

$ws.Range("H137").Value = 9012.200000000001
$ws.Range("J137").Value = 13343.667
$ws.Range("L137").Value = 40031.001
$ws.Range("N137").Value = -50231.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5074
$ws.Range("I82").Value = 8265
$ws.Range("J82").Value = 3159.4
$ws.Range("K82").Value = 8265
$ws.Range("L82").Value = 3159.4
$ws.Range("M82").Value = -7904
$ws.Range("N82").Value = -3881.4

$ws.Range("H85").Value = 5074
$ws.Range("I85").Value = 8265
$ws.Range("J85").Value = 3159.4
$ws.Range("K85").Value = 8265
$ws.Range("L85").Value = 3159.4
$ws.Range("M85").Value = -7017
$ws.Range("N85").Value = -5655.4

$ws.Range("H136").Value = 9349.875
$ws.Range("I136").Value = 3266.3333
$ws.Range("J136").Value = 13000
$ws.Range("K136").Value = 9798.999899999999
$ws.Range("L136").Value = 39000
$ws.Range("M136").Value = -7248.999899999999
$ws.Range("N136").Value = -44100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 42962.855
$ws.Range("J81").Value = 9374.5
$ws.Range("L81").Value = 18749
$ws.Range("N81").Value = -20871

$ws.Range("H84").Value = 42962.855
$ws.Range("J84").Value = 9374.5
$ws.Range("L84").Value = 93745
$ws.Range("N84").Value = -104353

$ws.Range("H136").Value = 5315
$ws.Range("I136").Value = 4794.4443
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 14383.3329
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -11833.3329
$ws.Range("N136").Value = -35100
